# Auto-generated edit script: update crypto price/volume table per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.289.04"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "2.229.16"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.67"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.92"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.403"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0900"
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "2.557.60"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.52"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.32"
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.799"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "2.232.17"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "42.168.29"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("D19").Value = "0.0₃0936"
$ws.Range("E19").Value = "  +4.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  +2.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.25"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "244.47"
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.68"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.86"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.140"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.34"
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("E30").Value = "  +2.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.67"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0652"
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.37"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("E39").Value = "  +6.10%  "
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.67"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("B42").Value = "TerraClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000225"
$ws.Range("E42").Value = "  -7.96%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0961"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.21"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.06"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "1.458.21"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.33"
$ws.Range("E47").Value = "  -9.61%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.08"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.07"
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("E51").Value = "  +3.96%  "
